# updated main GSC export data
# Appends the 2025-12-25 row to the "Chart" sheet (mirrors the upstream
# GSC export growing by one more day), matching the pattern of all the
# prior date rows: Date | Non-HTTPS URLs | HTTPS URLs.

$wb = $excel.ActiveWorkbook

$chart = $wb.Worksheets.Item(1)   # "Chart" sheet
$table = $wb.Worksheets.Item(2)   # "Table" sheet (unaffected in content)

$newRow = 81

# Write the date as a literal formula-derived string first and then paste
# it back as a value. This keeps the cell a genuine text/shared-string
# cell (matching every other date cell in the column) instead of letting
# the date-like text get auto-converted into a real date serial number,
# and it does so without adding a stray unused cell style.
$dateCell = $chart.Range("A" + $newRow)
$dateCell.Formula = "=""2025-12-25"""
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)  # xlPasteValues

$chart.Range("B" + $newRow).Value = 0
$chart.Range("C" + $newRow).Value = 31
